$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 419, shifting the existing data (rows 419-451) down to 421-453.
$ws.Rows.Item(419).Resize(2).Insert()

# New row 419: Coliflor "Primera" entry for the latest week (date serial 45132).
$ws.Range("A419").Value = 11
$ws.Range("B419").Value = "Vega Monumental Concepción"
$ws.Range("C419").Value = "Bíobío"
$ws.Range("D419").Value = 45132
$ws.Range("E419").Value = 8
$ws.Range("F419").Value = 100112008
$ws.Range("G419").Value = "Coliflor"
$ws.Range("H419").Value = "Sin especificar"
$ws.Range("I419").Value = "Primera"
$ws.Range("J419").Value = 2000
$ws.Range("K419").Value = 900
$ws.Range("L419").Value = 1000
$ws.Range("M419").Value = 950
$ws.Range("N419").Value = "`$/unidad"
$ws.Range("O419").Value = "Región Metropolitana"
$ws.Range("P419").Value = 950
$ws.Range("Q419").Value = 1
$ws.Range("R419").Value = "Hortaliza"

# New row 420: Coliflor "Segunda" entry for the latest week (date serial 45132).
$ws.Range("A420").Value = 11
$ws.Range("B420").Value = "Vega Monumental Concepción"
$ws.Range("C420").Value = "Bíobío"
$ws.Range("D420").Value = 45132
$ws.Range("E420").Value = 8
$ws.Range("F420").Value = 100112008
$ws.Range("G420").Value = "Coliflor"
$ws.Range("H420").Value = "Sin especificar"
$ws.Range("I420").Value = "Segunda"
$ws.Range("J420").Value = 1000
$ws.Range("K420").Value = 700
$ws.Range("L420").Value = 700
$ws.Range("M420").Value = 700
$ws.Range("N420").Value = "`$/unidad"
$ws.Range("O420").Value = "Región Metropolitana"
$ws.Range("P420").Value = 700
$ws.Range("Q420").Value = 1
$ws.Range("R420").Value = "Hortaliza"
